# Apply updated crypto price/volume data (GitHub Actions scrape refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.107.38'
$ws.Range('E2').Value = '  -0.35%  '

$ws.Range('D3').Value = '1.653.00'
$ws.Range('E3').Value = '  -0.40%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '218.49'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.12%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.5291'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +1.44%  '

$ws.Range('E8').Value = '  -2.21%  '

$ws.Range('E9').Value = '  +0.03%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '20.41'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -2.92%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07743'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +0.41%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '4.482'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +1.31%  '

$ws.Range('D13').Value = '1.657.01'
$ws.Range('E13').Value = '  -0.89%  '

$ws.Range('E14').Value = '  +0.00%  '

$ws.Range('D15').Value = '0.0₅8131'

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '65.23'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +0.69%  '

$ws.Range('D17').Value = '26.126.29'
$ws.Range('E17').Value = '  -0.41%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '1.003'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -0.31%  '

$ws.Range('E19').Value = '  -2.61%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '194.18'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +0.97%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '10.04'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -0.95%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '5.985'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -1.52%  '

$ws.Range('E23').Value = '  -0.39%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '140.05'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +1.20%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.1242'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +0.35%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '7.262'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +0.64%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '16.17'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +0.29%  '

$ws.Range('E28').Value = '  +1.55%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.05934'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -1.01%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.278'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -0.18%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '3.504'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -3.42%  '

$ws.Range('E32').Value = '  -2.49%  '

$ws.Range('E33').Value = '  -5.57%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '2.412'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -0.03%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.9440'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -3.46%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.758'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -0.79%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.5628'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -4.27%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.01604'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +1.16%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '5.851'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -1.50%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.8453'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -2.23%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.003'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -0.17%  '

$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D42').Value = '1.009.29'
$ws.Range('E42').Value = '  -2.16%  '

$ws.Range('B43').Value = 'Quant'
$ws.Range('C43').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '100.87'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +1.25%  '

$ws.Range('D44').Value = '1.798.85'
$ws.Range('E44').Value = '  -0.19%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '56.82'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -0.33%  '

$ws.Range('E46').Value = '  -3.69%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.003'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -0.63%  '

$ws.Range('E48').Value = '  +1.36%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.05150'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -0.68%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.470'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +0.61%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '7.756'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -4.24%  '
